$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5 (Kerosene row) to make room for
# the "Aviation Gasoline" entry that belongs right after Motor Gasoline.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "Aviation Gasoline"
$ws.Range("B5").Value = "07_petroleum_products"
$ws.Range("C5").Value = "07_02_aviation_gasoline"
